$d = $word.ActiveDocument

$oldText = "Datas das campanhas de 2018 que usam Perseu: 30 de outubro a 8 de novembro e 29 de novembro a 8 de dezembro"
$newText = "Datas das campanhas de Bootes: 14-23 de maio, 13-22 de junho, 12-21 de julho"

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs($i)
    $r = $para.Range
    # Exclude the trailing paragraph mark from the comparison/replacement.
    $cmp = $r.Duplicate
    [void]$cmp.MoveEnd(1, -1)
    if ($cmp.Text -eq $oldText) {
        $cmp.Text = ""
        $cmp.InsertAfter($newText)
    }
}
